$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.256.97'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.87%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.719.23'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.35%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4713'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2624'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06206'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.718.12'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07084'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.34'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.5990'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.432'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.22'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.66%  '
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.269.71'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.94%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006812'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.54'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.937.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.539'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.731'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.295'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '134.55'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('E26').Value = '  +0.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.402'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.764'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.49%  '
$ws.Range('E29').Value = '  +1.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.977'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.677'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.07762'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04467'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.617'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9763'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6182'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9348'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '113.14'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +17.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.422'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.921'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.000'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E42').Value = '  +0.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.484'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +13.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3828'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1183'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.283'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05272'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.784'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.26'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3382'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.217'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.63%  '
